$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44607
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 11000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 11500
$ws.Range("S2").Value = 639

# Row 3
$ws.Range("D3").Value = 44607
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 240
$ws.Range("N3").Value = 9000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 9500
$ws.Range("S3").Value = 528

# Row 4
$ws.Range("D4").Value = 44606
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 240
$ws.Range("N4").Value = 11500
$ws.Range("P4").Value = 11750
$ws.Range("S4").Value = 653

# Row 5
$ws.Range("D5").Value = 44606
$ws.Range("L5").Value = "Segunda"
$ws.Range("N5").Value = 9500
$ws.Range("O5").Value = 10000
$ws.Range("P5").Value = 9750
$ws.Range("S5").Value = 542

# Row 6
$ws.Range("D6").Value = 44610
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 13000
$ws.Range("O6").Value = 14000
$ws.Range("P6").Value = 13500
$ws.Range("S6").Value = 750

# Row 7
$ws.Range("D7").Value = 44610
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 200

# Row 8
$ws.Range("D8").Value = 44687
$ws.Range("L8").Value = "Especial"
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 19000
$ws.Range("P8").Value = 18500
$ws.Range("S8").Value = 1028

# Row 9
$ws.Range("D9").Value = 44687
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 14500
$ws.Range("S9").Value = 806

# Row 10
$ws.Range("D10").Value = 44595
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 15500
$ws.Range("P10").Value = 15750
$ws.Range("S10").Value = 875

# Row 11
$ws.Range("D11").Value = 44634
$ws.Range("M11").Value = 200

# Row 12
$ws.Range("D12").Value = 44634

# Row 16
$ws.Range("D16").Value = 44630
$ws.Range("M16").Value = 300
$ws.Range("N16").Value = 15000
$ws.Range("O16").Value = 16000
$ws.Range("P16").Value = 15500
$ws.Range("S16").Value = 861

# Row 17
$ws.Range("D17").Value = 44630
$ws.Range("M17").Value = 300
$ws.Range("N17").Value = 12000
$ws.Range("O17").Value = 13000
$ws.Range("P17").Value = 12500
$ws.Range("S17").Value = 694

# Row 18
$ws.Range("D18").Value = 44630
$ws.Range("L18").Value = "Segunda"
$ws.Range("N18").Value = 9000
$ws.Range("O18").Value = 10000
$ws.Range("P18").Value = 9500
$ws.Range("S18").Value = 528

# Row 19
$ws.Range("D19").Value = 44603
$ws.Range("L19").Value = "Especial"
$ws.Range("M19").Value = 240
$ws.Range("N19").Value = 14500
$ws.Range("O19").Value = 15000
$ws.Range("P19").Value = 14750
$ws.Range("S19").Value = 819

# Row 20
$ws.Range("D20").Value = 44685
$ws.Range("L20").Value = "Especial"
$ws.Range("N20").Value = 19000
$ws.Range("O20").Value = 20000
$ws.Range("P20").Value = 19500
$ws.Range("S20").Value = 1083

# Row 21
$ws.Range("D21").Value = 44685
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 160

# Row 22
$ws.Range("D22").Value = 44294
$ws.Range("L22").Value = "Especial"
$ws.Range("M22").Value = 200
$ws.Range("N22").Value = 14500
$ws.Range("O22").Value = 15000
$ws.Range("P22").Value = 14750
$ws.Range("S22").Value = 819

# Row 23
$ws.Range("D23").Value = 44294
$ws.Range("L23").Value = "Primera"
$ws.Range("N23").Value = 12500
$ws.Range("O23").Value = 13000
$ws.Range("P23").Value = 12750
$ws.Range("S23").Value = 708

# Row 24
$ws.Range("L24").Value = "Segunda"
$ws.Range("M24").Value = 240
$ws.Range("N24").Value = 10500
$ws.Range("O24").Value = 11000
$ws.Range("P24").Value = 10750
$ws.Range("S24").Value = 597

# Row 25
$ws.Range("D25").Value = 44636
$ws.Range("L25").Value = "Especial"
$ws.Range("N25").Value = 14000
$ws.Range("O25").Value = 15000
$ws.Range("P25").Value = 14500
$ws.Range("S25").Value = 806

# Row 26
$ws.Range("D26").Value = 44636
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 200
$ws.Range("N26").Value = 10000
$ws.Range("P26").Value = 10500
$ws.Range("S26").Value = 583

# Row 27
$ws.Range("D27").Value = 44609
$ws.Range("L27").Value = "Primera"
$ws.Range("N27").Value = 13000
$ws.Range("O27").Value = 14000
$ws.Range("P27").Value = 13500
$ws.Range("S27").Value = 750

# Row 28
$ws.Range("D28").Value = 44609
$ws.Range("L28").Value = "Segunda"
$ws.Range("M28").Value = 240
$ws.Range("N28").Value = 11000
$ws.Range("O28").Value = 12000
$ws.Range("P28").Value = 11500
$ws.Range("S28").Value = 639

# Row 29
$ws.Range("D29").Value = 44631
$ws.Range("L29").Value = "Especial"
$ws.Range("M29").Value = 240
$ws.Range("N29").Value = 15000
$ws.Range("O29").Value = 16000
$ws.Range("P29").Value = 15500
$ws.Range("S29").Value = 861

# Row 30
$ws.Range("D30").Value = 44631
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 248
$ws.Range("N30").Value = 12000
$ws.Range("O30").Value = 13000
$ws.Range("P30").Value = 12516
$ws.Range("S30").Value = 695

# Row 31
$ws.Range("D31").Value = 44631
$ws.Range("L31").Value = "Segunda"
$ws.Range("N31").Value = 9000
$ws.Range("O31").Value = 10000
$ws.Range("P31").Value = 9500
$ws.Range("S31").Value = 528

# Row 32
$ws.Range("D32").Value = 44637
$ws.Range("L32").Value = "Especial"
$ws.Range("M32").Value = 200
$ws.Range("N32").Value = 14000
$ws.Range("O32").Value = 15000
$ws.Range("P32").Value = 14500
$ws.Range("S32").Value = 806

# Row 33
$ws.Range("D33").Value = 44637
$ws.Range("L33").Value = "Primera"
$ws.Range("N33").Value = 10000
$ws.Range("O33").Value = 11000
$ws.Range("P33").Value = 10500
$ws.Range("S33").Value = 583
